$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.072.48"
$ws.Range("E2").Value = "  -0.21%  "

# Row 3
$ws.Range("D3").Value = "1.636.14"
$ws.Range("E3").Value = "  -1.89%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").Value = "'213.50"
$ws.Range("E5").Value = "  +1.41%  "

# Row 6
$ws.Range("D6").Value = "'0.5241"
$ws.Range("E6").Value = "  -0.23%  "

# Row 7
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("D8").Value = "'0.2596"
$ws.Range("E8").Value = "  -1.24%  "

# Row 9
$ws.Range("E9").Value = "  +0.22%  "

# Row 10
$ws.Range("E10").Value = "  -2.29%  "

# Row 11
$ws.Range("D11").Value = "'0.07656"
$ws.Range("E11").Value = "  +1.57%  "

# Row 12
$ws.Range("D12").Value = "1.633.35"
$ws.Range("E12").Value = "  -2.47%  "

# Row 13
$ws.Range("D13").Value = "'4.415"
$ws.Range("E13").Value = "  -0.62%  "

# Row 14
$ws.Range("D14").Value = "1.859.71"
$ws.Range("E14").Value = "  -1.99%  "

# Row 15
$ws.Range("D15").Value = "'0.5488"
$ws.Range("E15").Value = "  -1.37%  "

# Row 16
$ws.Range("D16").Value = "0.0₅8192"
$ws.Range("E16").Value = "  +3.14%  "

# Row 17
$ws.Range("D17").Value = "'64.99"
$ws.Range("E17").Value = "  -2.68%  "

# Row 18
$ws.Range("D18").Value = "26.065.81"
$ws.Range("E18").Value = "  -0.39%  "

# Row 20
$ws.Range("D20").Value = "'4.685"
$ws.Range("E20").Value = "  -1.30%  "

# Row 21
$ws.Range("D21").Value = "'188.37"
$ws.Range("E21").Value = "  +1.02%  "

# Row 22
$ws.Range("E22").Value = "  -2.16%  "

# Row 23
$ws.Range("D23").Value = "'6.150"
$ws.Range("E23").Value = "  -0.48%  "

# Row 24
$ws.Range("E24").Value = "  +0.00%  "

# Row 25
$ws.Range("D25").Value = "'145.62"
$ws.Range("E25").Value = "  -2.75%  "

# Row 26
$ws.Range("D26").Value = "'0.1213"
$ws.Range("E26").Value = "  -3.07%  "

# Row 27
$ws.Range("D27").Value = "'7.395"
$ws.Range("E27").Value = "  -1.49%  "

# Row 28
$ws.Range("D28").Value = "'15.79"
$ws.Range("E28").Value = "  -1.04%  "

# Row 29
$ws.Range("D29").Value = "'1.402"

# Row 30
$ws.Range("D30").Value = "'0.05964"
$ws.Range("E30").Value = "  -4.39%  "

# Row 31
$ws.Range("D31").Value = "'1.254"
$ws.Range("E31").Value = "  -2.01%  "

# Row 32
$ws.Range("D32").Value = "'3.434"
$ws.Range("E32").Value = "  -2.24%  "

# Row 33
$ws.Range("D33").Value = "'3.402"
$ws.Range("E33").Value = "  -0.50%  "

# Row 34
$ws.Range("D34").Value = "'1.637"
$ws.Range("E34").Value = "  +0.40%  "

# Row 35
$ws.Range("D35").Value = "'0.9862"
$ws.Range("E35").Value = "  -1.17%  "

# Row 36
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.396"
$ws.Range("E36").Value = "  -0.69%  "

# Row 37
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'2.762"
$ws.Range("E37").Value = "  +1.15%  "

# Row 38
$ws.Range("D38").Value = "'0.5713"
$ws.Range("E38").Value = "  -5.56%  "

# Row 39
$ws.Range("D39").Value = "'0.01616"
$ws.Range("E39").Value = "  +0.09%  "

# Row 40
$ws.Range("D40").Value = "'0.8549"
$ws.Range("E40").Value = "  -1.84%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'5.737"
$ws.Range("E41").Value = "  -6.51%  "

# Row 42
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  -0.21%  "

# Row 43
$ws.Range("D43").Value = "1.035.90"
$ws.Range("E43").Value = "  -6.27%  "

# Row 44
$ws.Range("D44").Value = "'100.57"
$ws.Range("E44").Value = "  +0.42%  "

# Row 45
$ws.Range("D45").Value = "1.786.74"
$ws.Range("E45").Value = "  -1.86%  "

# Row 46
$ws.Range("D46").Value = "0.0₈107"
$ws.Range("E46").Value = "  +4.02%  "

# Row 47
$ws.Range("D47").Value = "'55.74"
$ws.Range("E47").Value = "  +0.46%  "

# Row 48
$ws.Range("D48").Value = "'1.000"
$ws.Range("E48").Value = "  -0.45%  "

# Row 49
$ws.Range("D49").Value = "'8.036"
$ws.Range("E49").Value = "  -0.28%  "

# Row 50
$ws.Range("E50").Value = "  -1.30%  "

# Row 51
$ws.Range("D51").Value = "'0.4223"
$ws.Range("E51").Value = "  -0.60%  "

Write-Host "Update complete"
